# Auto-generated edit script: apply updated cell values per diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 600
$ws.Range("J2").Value = 750
$ws.Range("L2").Value = 750
$ws.Range("N2").Value = -976
$ws.Range("H5").Value = 356.66666
$ws.Range("I5").Value = 313.75
$ws.Range("J5").Value = 700
$ws.Range("K5").Value = 313.75
$ws.Range("L5").Value = 700
$ws.Range("M5").Value = -198.75
$ws.Range("N5").Value = -930
$ws.Range("H11").Value = 20000208
$ws.Range("I11").Value = 20000208
$ws.Range("K11").Value = 20000208
$ws.Range("M11").Value = -20000068
$ws.Range("H18").Value = 17271.273
$ws.Range("I18").Value = 4500
$ws.Range("J18").Value = 20109.334
$ws.Range("K18").Value = 4500
$ws.Range("L18").Value = 20109.334
$ws.Range("M18").Value = -4216
$ws.Range("N18").Value = -20677.334
$ws.Range("H70").Value = 12183.777
$ws.Range("I70").Value = 1050
$ws.Range("J70").Value = 15364.857
$ws.Range("K70").Value = 3150
$ws.Range("L70").Value = 46094.571
$ws.Range("M70").Value = -2880
$ws.Range("N70").Value = -46634.571
$ws.Range("H73").Value = 12183.777
$ws.Range("I73").Value = 1050
$ws.Range("J73").Value = 15364.857
$ws.Range("K73").Value = 3150
$ws.Range("L73").Value = 46094.571
$ws.Range("M73").Value = -2214
$ws.Range("N73").Value = -47966.571
$ws.Range("H100").Value = 1742.4166
$ws.Range("I100").Value = 1719
$ws.Range("K100").Value = 1719
$ws.Range("M100").Value = -1178
$ws.Range("H107").Value = 517.4
$ws.Range("I107").Value = 436.25
$ws.Range("K107").Value = 436.25
$ws.Range("M107").Value = 1483.75
$ws.Range("H113").Value = 24828.223
$ws.Range("I113").Value = 24828.223
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 24828.223
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -21574.223
$ws.Range("N113").ClearContents()
$ws.Range("H138").Value = 2872.4524
$ws.Range("J138").Value = 3619.5
$ws.Range("L138").Value = 10858.5
$ws.Range("N138").Value = -21138.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 382.85715
$ws.Range("I5").Value = 222
$ws.Range("K5").Value = 222
$ws.Range("M5").Value = -110
$ws.Range("H32").Value = 2774.5059
$ws.Range("I32").Value = 2131.5916
$ws.Range("K32").Value = 2131.5916
$ws.Range("M32").Value = -1844.5916
$ws.Range("H74").Value = 1196.0513
$ws.Range("I74").Value = 967.10345
$ws.Range("K74").Value = 967.10345
$ws.Range("M74").Value = -93.10344999999995
$ws.Range("H77").Value = 1196.0513
$ws.Range("I77").Value = 967.10345
$ws.Range("K77").Value = 4835.51725
$ws.Range("M77").Value = -467.5172499999999
$ws.Range("H109").Value = 56437.5
$ws.Range("J109").Value = 56437.5
$ws.Range("L109").Value = 56437.5
$ws.Range("N109").Value = -59211.5
$ws.Range("H110").Value = 3185.1667
$ws.Range("I110").Value = 1819.6
$ws.Range("J110").Value = 10013
$ws.Range("K110").Value = 1819.6
$ws.Range("L110").Value = 10013
$ws.Range("M110").Value = 225.4000000000001
$ws.Range("N110").Value = -14103
$ws.Range("H122").Value = 3107.1667
$ws.Range("I122").Value = 1728.6
$ws.Range("K122").Value = 5185.799999999999
$ws.Range("M122").Value = -2735.799999999999
$ws.Range("H139").Value = 28451.4
$ws.Range("J139").Value = 32064.25
$ws.Range("L139").Value = 32064.25
$ws.Range("N139").Value = -42344.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 382.85715
$ws.Range("I4").Value = 222
$ws.Range("K4").Value = 222
$ws.Range("M4").Value = -107
$ws.Range("H22").Value = 899.6667
$ws.Range("I22").Value = 849.5
$ws.Range("K22").Value = 849.5
$ws.Range("M22").Value = -676.5
$ws.Range("H97").Value = 17500
$ws.Range("I97").Value = 10000
$ws.Range("J97").Value = 20000
$ws.Range("K97").Value = 10000
$ws.Range("L97").Value = 20000
$ws.Range("M97").Value = -9009
$ws.Range("N97").Value = -21982
$ws.Range("H107").Value = 1475.75
$ws.Range("I107").Value = 1369.421
$ws.Range("K107").Value = 1369.421
$ws.Range("M107").Value = 550.579

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 867.4545000000001
$ws.Range("J16").Value = 998.6667
$ws.Range("L16").Value = 998.6667
$ws.Range("N16").Value = -1572.6667
$ws.Range("H22").Value = 1062.25
$ws.Range("J22").Value = 1283
$ws.Range("L22").Value = 1283
$ws.Range("N22").Value = -1983
$ws.Range("H86").Value = 2058.4
$ws.Range("I86").Value = 1640.8572
$ws.Range("K86").Value = 1640.8572
$ws.Range("M86").Value = -517.8571999999999
$ws.Range("H89").Value = 2058.4
$ws.Range("I89").Value = 1640.8572
$ws.Range("K89").Value = 8204.286
$ws.Range("M89").Value = -2588.286
$ws.Range("H113").Value = 867.4545000000001
$ws.Range("J113").Value = 998.6667
$ws.Range("L113").Value = 998.6667
$ws.Range("N113").Value = -5338.6667
$ws.Range("H122").Value = 3079.1052
$ws.Range("I122").Value = 1920.4
$ws.Range("J122").Value = 7424.25
$ws.Range("K122").Value = 5761.200000000001
$ws.Range("L122").Value = 22272.75
$ws.Range("M122").Value = -3311.200000000001
$ws.Range("N122").Value = -27172.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 909.73334
$ws.Range("I107").Value = 300
$ws.Range("J107").Value = 953.2857
$ws.Range("K107").Value = 900
$ws.Range("L107").Value = 2859.8571
$ws.Range("M107").Value = 1020
$ws.Range("N107").Value = -6699.8571
$ws.Range("H122").Value = 992.8182
$ws.Range("J122").Value = 1225.8572
$ws.Range("L122").Value = 11032.7148
$ws.Range("N122").Value = -15932.7148
$ws.Range("H131").Value = 800.61
$ws.Range("J131").Value = 806.81445
$ws.Range("L131").Value = 2420.44335
$ws.Range("N131").Value = -12500.44335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3888
$ws.Range("J80").Value = 3999
$ws.Range("L80").Value = 3999
$ws.Range("N80").Value = -5995
$ws.Range("H83").Value = 3888
$ws.Range("J83").Value = 3999
$ws.Range("L83").Value = 19995
$ws.Range("N83").Value = -29979
$ws.Range("H92").Value = 19423.625
$ws.Range("J92").Value = 19423.625
$ws.Range("L92").Value = 19423.625
$ws.Range("N92").Value = -23167.625
$ws.Range("H98").Value = 17661
$ws.Range("J98").Value = 17661
$ws.Range("L98").Value = 17661
$ws.Range("N98").Value = -23651
$ws.Range("H99").Value = 4500
$ws.Range("I99").Value = 4500
$ws.Range("K99").Value = 4500
$ws.Range("M99").Value = -2254
$ws.Range("H102").Value = 2584.2173
$ws.Range("I102").Value = 2607.611
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 2607.611
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = -985.6109999999999
$ws.Range("N102").Value = -5744

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 5000
$ws.Range("I35").Value = 5000
$ws.Range("K35").Value = 5000
$ws.Range("M35").Value = -4664
$ws.Range("H46").Value = 2298.3845
$ws.Range("I46").Value = 1220
$ws.Range("J46").Value = 2972.375
$ws.Range("K46").Value = 1220
$ws.Range("L46").Value = 2972.375
$ws.Range("M46").Value = -1032
$ws.Range("N46").Value = -3348.375
$ws.Range("H61").Value = 2224
$ws.Range("I61").Value = 2187.7778
$ws.Range("J61").Value = 2332.6667
$ws.Range("K61").Value = 2187.7778
$ws.Range("L61").Value = 2332.6667
$ws.Range("M61").Value = -1985.7778
$ws.Range("N61").Value = -2736.6667
$ws.Range("H106").Value = 19948
$ws.Range("J106").Value = 19948
$ws.Range("L106").Value = 19948
$ws.Range("N106").Value = -22472
$ws.Range("H113").Value = 2224
$ws.Range("I113").Value = 2187.7778
$ws.Range("J113").Value = 2332.6667
$ws.Range("K113").Value = 2187.7778
$ws.Range("L113").Value = 2332.6667
$ws.Range("M113").Value = -17.77779999999984
$ws.Range("N113").Value = -6672.6667
$ws.Range("H122").Value = 6032.5
$ws.Range("I122").Value = 1548.75
$ws.Range("J122").Value = 15000
$ws.Range("K122").Value = 4646.25
$ws.Range("L122").Value = 45000
$ws.Range("M122").Value = -2196.25
$ws.Range("N122").Value = -49900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 47343.5
$ws.Range("J95").Value = 47343.5
$ws.Range("L95").Value = 47343.5
$ws.Range("N95").Value = -52835.5
$ws.Range("H113").Value = 628.9474
$ws.Range("I113").Value = 423
$ws.Range("K113").Value = 1269
$ws.Range("M113").Value = 901
$ws.Range("H126").Value = 13409.272
$ws.Range("I126").Value = 16571.715
$ws.Range("K126").Value = 49715.145
$ws.Range("M126").Value = -47245.145
$ws.Range("H132").Value = 1048.3715
$ws.Range("I132").Value = 828.8387
$ws.Range("K132").Value = 2486.5161
$ws.Range("M132").Value = 43.48390000000018
